## Handles float input without breaking stuff
## -------------------------------------------------------------
## Fixes the marksheet summary block (rows 10-12) so the "Marking"
## row's negative-marking value is stored as a real number instead
## of text (the "float input" that used to break downstream math),
## refreshes the computed totals, and rebuilds the per-question
## answer-key block (rows 15-40) so the student's answer is shown in
## column A (colour-coded correct/incorrect) next to the existing
## "Correct Ans" column B. The redundant second/third Student-Ans /
## Correct-Ans column pairs (D:E beyond row 18, and G:H entirely)
## are removed, which also shrinks the sheet's used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Summary block (rows 10-12): give the row labels the header
#    style and correct the Right/Wrong/Not-Attempt/Max numbers.
#    C11 especially moves from a text "-1" to a real numeric -1.
# ---------------------------------------------------------------
$ws.Range("A10").Style = "normalStyle"
$ws.Range("A11").Style = "normalStyle"
$ws.Range("A12").Style = "normalStyle"

$ws.Range("B10").Value = 25
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 2
$ws.Range("E10").Value = 28

$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("B12").Value = 100
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "99/112"

# ---------------------------------------------------------------
# 2. Per-question answer rows (15-40): column A now carries the
#    student's recorded answer, styled green ("correctStyle") when
#    it matches the "Correct Ans" in column B, red
#    ("incorrectStyle") when it does not, or left with the plain
#    ("normalStyle") look for a couple of legacy rows. Rows 16-18
#    additionally keep a second Student/Correct pair in D:E.
# ---------------------------------------------------------------
$answers = @(
    @{Row=16; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=16; Col="D"; Text="Option A"; Style="correctStyle"},
    @{Row=17; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=17; Col="D"; Text="Option C"; Style="correctStyle"},
    @{Row=18; Col="A"; Text="Option B"; Style="correctStyle"},
    @{Row=18; Col="D"; Text="Option D"; Style="correctStyle"},
    @{Row=19; Col="A"; Text="Option C"; Style="correctStyle"},
    @{Row=20; Col="A"; Text="Option B"; Style="correctStyle"},
    @{Row=21; Col="A"; Text="Option C"; Style="correctStyle"},
    @{Row=22; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=23; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=24; Col="A"; Text="Option A"; Style="normalStyle"},
    @{Row=25; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=26; Col="A"; Text="Option D"; Style="incorrectStyle"},
    @{Row=27; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=28; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=29; Col="A"; Text="Option D"; Style="normalStyle"},
    @{Row=30; Col="A"; Text="Option B"; Style="correctStyle"},
    @{Row=31; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=32; Col="A"; Text="Option C"; Style="correctStyle"},
    @{Row=33; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=34; Col="A"; Text="Option B"; Style="correctStyle"},
    @{Row=35; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=36; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=37; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=38; Col="A"; Text="Option A"; Style="correctStyle"},
    @{Row=39; Col="A"; Text="Option D"; Style="correctStyle"},
    @{Row=40; Col="A"; Text="Option D"; Style="correctStyle"}
)

foreach ($item in $answers) {
    $addr = "$($item.Col)$($item.Row)"
    $ws.Range($addr).Value = $item.Text
    $ws.Range($addr).Style = $item.Style
}

# ---------------------------------------------------------------
# 3. Drop the redundant Student Ans / Correct Ans column pairs:
#    D:E stay only for rows 16-18, so wipe them for rows 19-40;
#    the G:H pair is not used anywhere any more, so remove those
#    whole columns (this also shrinks the sheet dimension down to
#    A5:E40, matching the trimmed layout).
# ---------------------------------------------------------------
$ws.Range("D19:E40").Clear()
$ws.Columns("G:H").Delete()
